$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($index, $xml) {
    $p = $d.Paragraphs($index)
    $null = $p.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) "Option A: Do nothing." (paragraph 10) -> green highlight
# ---------------------------------------------------------------------------
$xml1 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:highlight w:val=`"green`"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t xml:space=`"preserve`">Option A: Do nothing. </w:t></w:r></w:p>"
Set-ParaXml 10 $xml1

# ---------------------------------------------------------------------------
# 2) "Option B: Add TTHM and HAA5 ..." (paragraph 14) -> green highlight + split
#    text with a new middle sentence inserted.
# ---------------------------------------------------------------------------
$xml2 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:highlight w:val=`"green`"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t xml:space=`"preserve`">Option B: Add TTHM and HAA5 relative risk tables to the appendix and reference them in the text. Edit certain text to explain why we combine these </w:t></w:r><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t xml:space=`"preserve`">(cumulative risk), what we understand about why these are separated into two MCLs, </w:t></w:r><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t xml:space=`"preserve`">and whether it makes a difference to separate them. Switch maps to 140 ug/l instead of 80.  </w:t></w:r></w:p>"
Set-ParaXml 14 $xml2

# ---------------------------------------------------------------------------
# 3-5) "Concentration ratios" / "Disparity ratios" / "Prevalence ratio"
#      (paragraphs 18, 19, 20) -> green highlight
# ---------------------------------------------------------------------------
$xml3 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:highlight w:val=`"green`"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>Concentration ratios</w:t></w:r></w:p>"
Set-ParaXml 18 $xml3

$xml4 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:highlight w:val=`"green`"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>Disparity ratios</w:t></w:r></w:p>"
Set-ParaXml 19 $xml4

$xml5 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:highlight w:val=`"green`"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>Prevalence ratio</w:t></w:r></w:p>"
Set-ParaXml 20 $xml5

# ---------------------------------------------------------------------------
# 6) "Option A: Sensitivity table with several alternative formulations"
#    (paragraph 26) -> green highlight
# ---------------------------------------------------------------------------
$xml6 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:highlight w:val=`"green`"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>Option A: Sensitivity table with several alternative formulations</w:t></w:r></w:p>"
Set-ParaXml 26 $xml6

# ---------------------------------------------------------------------------
# 7) "Option B: Do nothing." (paragraph 27, PFAS section) -> green highlight
#    plus "Add footnote" appended as new runs.
# ---------------------------------------------------------------------------
$xml7 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:highlight w:val=`"green`"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t xml:space=`"preserve`">Option B: Do nothing. </w:t></w:r><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t xml:space=`"preserve`">Add </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>footnote</w:t></w:r><w:proofErr w:type=`"gramEnd`"/></w:p>"
Set-ParaXml 27 $xml7

# ---------------------------------------------------------------------------
# 8) Next-Steps list reshuffle:
#    33: "Move all in-text edits over to the Overleaf file"
#        -> "Double-checking code is updated "
#    34: "Edit Overleaf file to address all non-controversial comment bubbles"
#        -> "Move all in-text edits over to the Overleaf file"
#    new paragraph inserted before 35 ("Incorporate latest H&M...")
#        -> "Edit Overleaf file to address all non-controversial comment bubbles"
# ---------------------------------------------------------------------------
$xml8 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">Double-checking code is </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>updated</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>"
Set-ParaXml 33 $xml8

$xml9 = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">Move all in-text edits over to the </w:t></w:r><w:r><w:t xml:space=`"preserve`">Overleaf </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>file</w:t></w:r><w:proofErr w:type=`"gramEnd`"/></w:p>"
Set-ParaXml 34 $xml9

$p34 = $d.Paragraphs(34)
$p34.Range.InsertParagraphAfter()
$xmlNew = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">Edit Overleaf file to address all non-controversial comment </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>bubbles</w:t></w:r><w:proofErr w:type=`"gramEnd`"/></w:p>"
Set-ParaXml 35 $xmlNew

# ---------------------------------------------------------------------------
# 9) New paragraph "Decide what to put in appendix" after "PFAS sensitivity
#    table" (now paragraph 43), followed by two extra blank paragraphs.
# ---------------------------------------------------------------------------
$pfas = $d.Paragraphs(43)
Write-Host "pfas check: $($pfas.Range.Text)"
$pfas.Range.InsertParagraphAfter()
$xmlAppendix = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">Decide </w:t></w:r><w:r><w:t xml:space=`"preserve`">what to put in </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>appendix</w:t></w:r><w:proofErr w:type=`"gramEnd`"/></w:p>"
Set-ParaXml 44 $xmlAppendix

$pAppendix = $d.Paragraphs(44)
$pAppendix.Range.InsertParagraphAfter()
$pAppendix.Range.InsertParagraphAfter()

Write-Host "edit complete"
